$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted before the old row 81, shifting
# every subsequent record down by one row (old row 186 becomes row 187).
$ws.Rows.Item(81).Insert()

$ws.Cells.Item(81, 1).Value = 6
$ws.Cells.Item(81, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(81, 3).Value = "Metropolitana"
$ws.Cells.Item(81, 4).Value = 44671
$ws.Cells.Item(81, 5).Value = 13
$ws.Cells.Item(81, 6).Value = 100112029
$ws.Cells.Item(81, 7).Value = "Orégano"
$ws.Cells.Item(81, 8).Value = "Sin especificar"
$ws.Cells.Item(81, 9).Value = "Primera"
$ws.Cells.Item(81, 10).Value = 31
$ws.Cells.Item(81, 11).Value = 12000
$ws.Cells.Item(81, 12).Value = 13000
$ws.Cells.Item(81, 13).Value = 12484
$ws.Cells.Item(81, 14).Value = "$/docena de atados"
$ws.Cells.Item(81, 15).Value = "Región Metropolitana"
$ws.Cells.Item(81, 16).Value = 4161
$ws.Cells.Item(81, 17).Value = 3
$ws.Cells.Item(81, 18).Value = "Hortaliza"

# Keep the date formatted like the other date cells in column D.
$ws.Cells.Item(81, 4).NumberFormat = $ws.Cells.Item(82, 4).NumberFormat
